# tokistorage-client-deck-en.pptx — remove Trial/Individual plans, keep
# Three-Generation and Timeless Transformation, add a self-service note.
# Slide 7 (1-based) = "INVESTMENT GUIDE" pricing slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# ---------------------------------------------------------------------
# 1. Sub-header text above the pricing cards.
# ---------------------------------------------------------------------
$s.Shapes.Item("TextBox 2").TextFrame.TextRange.Text = `
    "A deeper proof of existence through dialogue. Timeless Consulting pricing."

# ---------------------------------------------------------------------
# 2. Turn the old "Trial Plan" card background into the slim self-service
#    banner (EFF6FF fill, no border, shorter height).
# ---------------------------------------------------------------------
$banner = $s.Shapes.Item("Rectangle 4")
$banner.Left = 36
$banner.Top = 82.8
$banner.Width = 619.2
$banner.Height = 39.6
$banner.Fill.ForeColor.RGB = RGBVal 0xEF 0xF6 0xFF
$banner.Line.Visible = $false

# ---------------------------------------------------------------------
# 3. New blue accent bar on the left edge of the banner.
# ---------------------------------------------------------------------
$accent = $s.Shapes.AddShape(1, 36, 82.8, 4.32, 39.6)
$accent.Name = "Rectangle 5"
$accent.Fill.ForeColor.RGB = RGBVal 0x25 0x63 0xEB
$accent.Line.Visible = $false

# ---------------------------------------------------------------------
# 4. New self-service note text inside the banner.
# ---------------------------------------------------------------------
$note = $s.Shapes.AddTextbox(1, 54, 86.4, 594, 32.4)
$note.Name = "TextBox 6"
$note.TextFrame.WordWrap = $true
$note.TextFrame.TextRange.Text = `
    "TokiQR (free), UV Laminate QR (¥5,000), and Quartz Glass (¥50,000) are also available as self-service."
$noteTr = $note.TextFrame.TextRange
$noteTr.Font.Size = 10
$noteTr.Font.Bold = $false
$noteTr.Font.Color.RGB = RGBVal 0x47 0x55 0x69
$noteTr.Font.Name = "Calibri"
$note.TextFrame.VerticalAnchor = 3   # msoAnchorMiddle
$note.TextFrame.WordWrap = $true
$note.TextFrame.TextRange.ParagraphFormat.Alignment = 1  # ppAlignLeft

# ---------------------------------------------------------------------
# 5. The old "Individual Plan" card background becomes the card behind
#    "Three-Generation" (re-positioned right under the banner).
# ---------------------------------------------------------------------
$threeGenCard = $s.Shapes.Item("Rectangle 9")
$threeGenCard.Left = 36
$threeGenCard.Top = 133.2
$threeGenCard.Width = 619.2
$threeGenCard.Height = 79.2

# ---------------------------------------------------------------------
# 6. The old "Trial" icon tile/letter/title/description get repurposed
#    as the "Three-Generation" icon tile/letter/title/description.
# ---------------------------------------------------------------------
$threeGenIconBg = $s.Shapes.Item("Rectangle 5")
$threeGenIconBg.Left = 46.8
$threeGenIconBg.Top = 144
$threeGenIconBg.Width = 36
$threeGenIconBg.Height = 36
$threeGenIconBg.Fill.ForeColor.RGB = RGBVal 0xFD 0xF8 0xE8

$threeGenIconLetter = $s.Shapes.Item("TextBox 6")
$threeGenIconLetter.Left = 46.8
$threeGenIconLetter.Top = 144
$threeGenIconLetter.Width = 36
$threeGenIconLetter.Height = 36
$threeGenIconLetter.TextFrame.TextRange.Text = "F"

$threeGenTitle = $s.Shapes.Item("TextBox 7")
$threeGenTitle.Left = 93.6
$threeGenTitle.Top = 141.84
$threeGenTitle.Width = 547.2
$threeGenTitle.Height = 23.04
$threeGenTitle.TextFrame.TextRange.Text = "Three-Generation — ¥550,000 ($3,600, recommended)"

$threeGenDesc = $s.Shapes.Item("TextBox 8")
$threeGenDesc.Left = 93.6
$threeGenDesc.Top = 166.32
$threeGenDesc.Width = 547.2
$threeGenDesc.Height = 37.44
$threeGenDesc.TextFrame.TextRange.Text = `
    "Three family members + on-site inscription on Sado Island. An alternative to traditional graves."

# ---------------------------------------------------------------------
# 7. Remove the now-unused "Individual Plan" pieces (icon tile, letter,
#    title, description) — their job is already done by the shapes above.
# ---------------------------------------------------------------------
$s.Shapes.Item("Rectangle 10").Delete()
$s.Shapes.Item("TextBox 11").Delete()
$s.Shapes.Item("TextBox 12").Delete()
$s.Shapes.Item("TextBox 13").Delete()

# ---------------------------------------------------------------------
# 8. Remove the now-unused old "Three-Generation" pieces (their old
#    positions/letters/colors were already reused above under the
#    banner, one row higher).
# ---------------------------------------------------------------------
$s.Shapes.Item("Rectangle 15").Delete()
$s.Shapes.Item("TextBox 16").Delete()
$s.Shapes.Item("TextBox 17").Delete()
$s.Shapes.Item("TextBox 18").Delete()

# ---------------------------------------------------------------------
# 9. The old "Three-Generation" card background becomes the card behind
#    "Timeless Transformation" (re-positioned to the second card slot).
# ---------------------------------------------------------------------
$timelessCard = $s.Shapes.Item("Rectangle 14")
$timelessCard.Left = 36
$timelessCard.Top = 223.2
$timelessCard.Width = 619.2
$timelessCard.Height = 79.2

# ---------------------------------------------------------------------
# 10. The old "Voice Memorial" icon tile/letter/title/description move
#     up into the second card slot and become "Timeless Transformation".
# ---------------------------------------------------------------------
$timelessIconBg = $s.Shapes.Item("Rectangle 20")
$timelessIconBg.Left = 46.8
$timelessIconBg.Top = 234
$timelessIconBg.Width = 36
$timelessIconBg.Height = 36

$timelessIconLetter = $s.Shapes.Item("TextBox 21")
$timelessIconLetter.Left = 46.8
$timelessIconLetter.Top = 234
$timelessIconLetter.Width = 36
$timelessIconLetter.Height = 36

$timelessTitle = $s.Shapes.Item("TextBox 22")
$timelessTitle.Left = 93.6
$timelessTitle.Top = 231.84
$timelessTitle.Width = 547.2
$timelessTitle.Height = 23.04
$timelessTitle.TextFrame.TextRange.Text = "Timeless Transformation — ¥5,550,000 ($36,500)"

$timelessDesc = $s.Shapes.Item("TextBox 23")
$timelessDesc.Left = 93.6
$timelessDesc.Top = 256.32
$timelessDesc.Width = 547.2
$timelessDesc.Height = 37.44
$timelessDesc.TextFrame.TextRange.Text = `
    "Fully bespoke. Voice restoration, custom design, and the complete premium experience."

# ---------------------------------------------------------------------
# 11. Remove the now-unused old "Voice Memorial" card background (its
#     job is already done by the repositioned card above).
# ---------------------------------------------------------------------
$s.Shapes.Item("Rectangle 19").Delete()

Write-Host "slide 7 pricing section rebuilt:" $s.Shapes.Count "shapes"
